$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.593.22'
$ws.Range("E2").Value = '  -2.82%  '
$ws.Range("D3").Value = '1.984.65'
$ws.Range("E3").Value = '  -3.75%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.639'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.74'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.21%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.04'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.361'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("E11").Value = '  -2.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.104'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.65%  '
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.79%  '
$ws.Range("D15").Value = '2.272.46'
$ws.Range("E15").Value = '  -3.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.60%  '
$ws.Range("D17").Value = '1.981.96'
$ws.Range("E17").Value = '  -3.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.70%  '
$ws.Range("D19").Value = '35.616.24'
$ws.Range("E19").Value = '  -2.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").Value = '0.0₃0846'
$ws.Range("E21").Value = '  -1.84%  '
$ws.Range("E22").Value = '  -1.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +22.03%  '
$ws.Range("E26").Value = '  -3.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.50%  '
$ws.Range("E30").Value = '  -2.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.31%  '
$ws.Range("E32").Value = '  -6.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0964'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +17.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0594'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.45'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.55%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  -3.68%  '
$ws.Range("E39").Value = '  +8.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.23'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.91'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.31%  '
$ws.Range("E42").Value = '  -2.02%  '
$ws.Range("E43").Value = '  -1.99%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.57%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0895'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.67%  '
$ws.Range("D48").Value = '1.374.76'
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("E49").Value = '  +1.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '46.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.52%  '
$ws.Range("E51").Value = '  -1.32%  '
